$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.729.00"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "1.657.85"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "303.09"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").Value = "0.3817"
$ws.Range("E7").Value = "  +0.93%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3620"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "51.27"
$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("D10").Value = "0.08192"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("D11").Value = "1.234"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("D13").Value = "22.59"
$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("D14").Value = "6.502"
$ws.Range("E14").Value = "  +0.81%  "

$ws.Range("D15").Value = "7.421"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").Value = "0.00001231"
$ws.Range("E16").Value = "  -0.55%  "

$ws.Range("D17").Value = "1.655.19"
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("D18").Value = "97.75"
$ws.Range("E18").Value = "  +2.42%  "

$ws.Range("D19").Value = "0.07004"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("D20").Value = "6.842"
$ws.Range("E20").Value = "  +3.81%  "

$ws.Range("D21").Value = "17.68"
$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "12.82"
$ws.Range("E23").Value = "  +2.60%  "

$ws.Range("D24").Value = "23.727.24"
$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("D25").Value = "2.514"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("E26").Value = "  -0.78%  "

$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("D28").Value = "152.63"
$ws.Range("E28").Value = "  +0.74%  "

$ws.Range("D29").Value = "5.220"
$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").Value = "134.31"
$ws.Range("E30").Value = "  +1.10%  "

$ws.Range("D31").Value = "1.837.13"
$ws.Range("E31").Value = "  +1.18%  "

$ws.Range("D32").Value = "7.043"
$ws.Range("E32").Value = "  +6.93%  "

$ws.Range("D33").Value = "2.218"
$ws.Range("E33").Value = "  +2.66%  "

$ws.Range("D34").Value = "11.97"
$ws.Range("E34").Value = "  +4.62%  "

$ws.Range("D35").Value = "1.062"
$ws.Range("E35").Value = "  -2.81%  "

$ws.Range("D36").Value = "0.02811"
$ws.Range("E36").Value = "  +1.91%  "

$ws.Range("D37").Value = "0.2520"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("D38").Value = "0.08807"
$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("D39").Value = "6.106"
$ws.Range("E39").Value = "  +2.06%  "

$ws.Range("D40").Value = "0.07045"
$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("D41").Value = "13.05"
$ws.Range("E41").Value = "  +7.34%  "

$ws.Range("D42").Value = "0.7023"
$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("D43").Value = "1.337"
$ws.Range("E43").Value = "  -1.18%  "

$ws.Range("D44").Value = "16.04"
$ws.Range("E44").Value = "  +3.12%  "

$ws.Range("D45").Value = "0.6527"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "2.315"
$ws.Range("E46").Value = "  +1.83%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("D48").Value = "3.964"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").Value = "0.07916"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").Value = "128.13"
$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("D51").Value = "1.185"
$ws.Range("E51").Value = "  -0.52%  "
